# Auto-generated edit script: refreshes market-price derived columns (H-N)
# across ALC/ARM/BSM/CRP/LTW/WVR sheets, per scheduled runner update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")  # row 31
$ws.Range("H31").Value = 2450
$ws.Range("J31").Value = 3050
$ws.Range("L31").Value = 9150
$ws.Range("N31").Value = -9610

$ws = $wb.Worksheets.Item("ALC")  # row 40
$ws.Range("H40").Value = 4145.8965
$ws.Range("I40").Value = 3957.0833
$ws.Range("K40").Value = 3957.0833
$ws.Range("M40").Value = -3782.0833

$ws = $wb.Worksheets.Item("ALC")  # row 132
$ws.Range("H132").Value = 1965.7931
$ws.Range("I132").Value = 1143.1786
$ws.Range("K132").Value = 3429.5358
$ws.Range("M132").Value = -899.5357999999997

$ws = $wb.Worksheets.Item("ALC")  # row 137
$ws.Range("H137").Value = 4331.24
$ws.Range("I137").Value = 1892
$ws.Range("K137").Value = 5676
$ws.Range("M137").Value = -3126

$ws = $wb.Worksheets.Item("ARM")  # row 32
$ws.Range("H32").Value = 14840.28
$ws.Range("I32").Value = 12121.833
$ws.Range("K32").Value = 12121.833
$ws.Range("M32").Value = -11834.833

$ws = $wb.Worksheets.Item("ARM")  # row 35
$ws.Range("H35").Value = 8081.8
$ws.Range("I35").Value = 8081.8
$ws.Range("K35").Value = 8081.8
$ws.Range("M35").Value = -7675.8

$ws = $wb.Worksheets.Item("ARM")  # row 45
$ws.Range("H45").Value = 5500
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 5500
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 5500
$ws.Range("M45").Value = ""
$ws.Range("N45").Value = -6254

$ws = $wb.Worksheets.Item("ARM")  # row 61
$ws.Range("H61").Value = 4183.5835
$ws.Range("I61").Value = 4209.364
$ws.Range("K61").Value = 4209.364
$ws.Range("M61").Value = -3997.364

$ws = $wb.Worksheets.Item("ARM")  # row 97
$ws.Range("H97").Value = 827.1429000000001
$ws.Range("I97").Value = 782.36365
$ws.Range("K97").Value = 782.36365
$ws.Range("M97").Value = -286.36365

$ws = $wb.Worksheets.Item("ARM")  # row 132
$ws.Range("H132").Value = 325
$ws.Range("I132").Value = 325
$ws.Range("K132").Value = 975
$ws.Range("M132").Value = 1555

$ws = $wb.Worksheets.Item("ARM")  # row 136
$ws.Range("H136").Value = 4183.5835
$ws.Range("I136").Value = 4209.364
$ws.Range("K136").Value = 12628.092
$ws.Range("M136").Value = -10078.092

$ws = $wb.Worksheets.Item("BSM")  # row 86
$ws.Range("H86").Value = 5028.88
$ws.Range("I86").Value = 3823.3572
$ws.Range("K86").Value = 3823.3572
$ws.Range("M86").Value = -2700.3572

$ws = $wb.Worksheets.Item("BSM")  # row 89
$ws.Range("H89").Value = 5028.88
$ws.Range("I89").Value = 3823.3572
$ws.Range("K89").Value = 19116.786
$ws.Range("M89").Value = -13500.786

$ws = $wb.Worksheets.Item("BSM")  # row 134
$ws.Range("H134").Value = 3141
$ws.Range("I134").Value = 1831.3334
$ws.Range("K134").Value = 5494.0002
$ws.Range("M134").Value = -2959.0002

$ws = $wb.Worksheets.Item("CRP")  # row 38
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = ""
$ws.Range("N38").Value = ""

$ws = $wb.Worksheets.Item("CRP")  # row 41
$ws.Range("H41").Value = 49507.5
$ws.Range("I41").Value = 250
$ws.Range("K41").Value = 250
$ws.Range("M41").Value = 178

$ws = $wb.Worksheets.Item("CRP")  # row 46
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = ""
$ws.Range("N46").Value = ""

$ws = $wb.Worksheets.Item("CRP")  # row 57
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").Value = ""

$ws = $wb.Worksheets.Item("CRP")  # row 59
$ws.Range("H59").Value = 30544.857
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 30544.857
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 30544.857
$ws.Range("M59").Value = ""
$ws.Range("N59").Value = -32834.857

$ws = $wb.Worksheets.Item("CRP")  # row 62
$ws.Range("H62").Value = 4912.8335
$ws.Range("I62").Value = 3500
$ws.Range("K62").Value = 3500
$ws.Range("M62").Value = -2876

$ws = $wb.Worksheets.Item("CRP")  # row 65
$ws.Range("H65").Value = 4912.8335
$ws.Range("I65").Value = 3500
$ws.Range("K65").Value = 17500
$ws.Range("M65").Value = -14380

$ws = $wb.Worksheets.Item("CRP")  # row 68
$ws.Range("H68").Value = 93759.86
$ws.Range("J68").Value = 93759.86
$ws.Range("L68").Value = 93759.86
$ws.Range("N68").Value = -95257.86

$ws = $wb.Worksheets.Item("CRP")  # row 71
$ws.Range("H71").Value = 93759.86
$ws.Range("J71").Value = 93759.86
$ws.Range("L71").Value = 281279.58
$ws.Range("N71").Value = -288767.58

$ws = $wb.Worksheets.Item("CRP")  # row 94
$ws.Range("H94").Value = 4069
$ws.Range("I94").Value = 1824.3334
$ws.Range("J94").Value = 6762.6
$ws.Range("K94").Value = 1824.3334
$ws.Range("L94").Value = 6762.6
$ws.Range("M94").Value = -1373.3334
$ws.Range("N94").Value = -7664.6

$ws = $wb.Worksheets.Item("LTW")  # row 7
$ws.Range("H7").Value = 5567.357
$ws.Range("I7").Value = 4720.3
$ws.Range("K7").Value = 4720.3
$ws.Range("M7").Value = -4608.3

$ws = $wb.Worksheets.Item("LTW")  # row 16
$ws.Range("H16").Value = 1474
$ws.Range("I16").Value = 965.6667
$ws.Range("J16").Value = 2999
$ws.Range("K16").Value = 965.6667
$ws.Range("L16").Value = 2999
$ws.Range("M16").Value = -795.6667
$ws.Range("N16").Value = -3339

$ws = $wb.Worksheets.Item("LTW")  # row 46
$ws.Range("H46").Value = 7141.467
$ws.Range("I46").Value = 5541
$ws.Range("J46").Value = 7541.5835
$ws.Range("K46").Value = 5541
$ws.Range("L46").Value = 7541.5835
$ws.Range("M46").Value = -5353
$ws.Range("N46").Value = -7917.5835

$ws = $wb.Worksheets.Item("LTW")  # row 126
$ws.Range("H126").Value = 5567.357
$ws.Range("I126").Value = 4720.3
$ws.Range("K126").Value = 14160.9
$ws.Range("M126").Value = -11690.9

$ws = $wb.Worksheets.Item("LTW")  # row 136
$ws.Range("H136").Value = 5230.5
$ws.Range("I136").Value = 5230.5
$ws.Range("K136").Value = 15691.5
$ws.Range("M136").Value = -13141.5

$ws = $wb.Worksheets.Item("WVR")  # row 54
$ws.Range("H54").Value = 62990.547
$ws.Range("J54").Value = 62990.547
$ws.Range("L54").Value = 62990.547
$ws.Range("N54").Value = -64030.547

$ws = $wb.Worksheets.Item("WVR")  # row 81
$ws.Range("H81").Value = 3388.3333
$ws.Range("J81").Value = 2499
$ws.Range("L81").Value = 4998
$ws.Range("N81").Value = -7120

$ws = $wb.Worksheets.Item("WVR")  # row 84
$ws.Range("H84").Value = 3388.3333
$ws.Range("J84").Value = 2499
$ws.Range("L84").Value = 24990
$ws.Range("N84").Value = -35598

$ws = $wb.Worksheets.Item("WVR")  # row 122
$ws.Range("H122").Value = 3173
$ws.Range("I122").Value = 2270
$ws.Range("K122").Value = 6810
$ws.Range("M122").Value = -4360

$ws = $wb.Worksheets.Item("WVR")  # row 126
$ws.Range("H126").Value = 6397
$ws.Range("I126").Value = 4030
$ws.Range("K126").Value = 12090
$ws.Range("M126").Value = -9620

$ws = $wb.Worksheets.Item("WVR")  # row 136
$ws.Range("H136").Value = 6237.375
$ws.Range("I136").Value = 5234.25
$ws.Range("K136").Value = 15702.75
$ws.Range("M136").Value = -13152.75
